# Update WTP Pain and Suffering values
# Values taken from tables 10 and 11 in doi: 10.1007/s10198-022-01512-3
# Dollar values are in AUD (2017)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Row 2: GI, severe (per day) ---
$ws.Range("C2").Value = 32
$ws.Range("D2").Value = 25
$ws.Range("E2").Value = 39

# --- Row 3: GI, mild (per day) ---
$ws.Range("C3").Value = 12
$ws.Range("D3").Value = 7
$ws.Range("E3").Value = 16

# --- Row 4: flulike, severe (per day) ---
$ws.Range("C4").Value = 23
$ws.Range("E4").Value = 30

# --- Row 5: flulike, mild (per day) ---
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = 13

# --- Row 6: IBS, severe (per year) - formula replaced with plain value ---
$ws.Range("C6").Value = 1367
$ws.Range("D6").Value = 938
$ws.Range("E6").Value = 1795

# --- Row 7: IBS, mild (per year) - formula replaced with plain value ---
$ws.Range("C7").Value = 575
$ws.Range("D7").Value = 272
$ws.Range("E7").Value = 877

# --- Row 8: ReA, severe (per year) - formula replaced with plain value ---
$ws.Range("C8").Value = 1412
$ws.Range("D8").Value = 915
$ws.Range("E8").Value = 1909

# --- Row 9: ReA, mild (per year) - formula replaced with plain value ---
$ws.Range("C9").Value = 531
$ws.Range("D9").Value = 156
$ws.Range("E9").Value = 907

# --- Row 10: HUS, severe (per year) - formula replaced with plain value ---
$ws.Range("C10").Value = 2195
$ws.Range("D10").Value = 1596
$ws.Range("E10").Value = 2794

# --- Row 11: HUS, mild (per year) - formula replaced with plain value ---
$ws.Range("C11").Value = 1025
$ws.Range("D11").Value = 630
$ws.Range("E11").Value = 1419

# --- Row 12: GBS, severe (per year) - formula replaced with plain value ---
$ws.Range("C12").Value = 1852
$ws.Range("D12").Value = 1337
$ws.Range("E12").Value = 2367

# --- Row 13: GBS, mild (per year) - formula replaced with plain value ---
$ws.Range("C13").Value = 835
$ws.Range("D13").Value = 453
$ws.Range("E13").Value = 1217

# Update the active selection to match the saved view state
$ws.Range("I13").Select()
